$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above the existing row 15 ("instructor_id", the bottom
# border of the "courses" table). Inserting shifts row 15 and everything
# below it down by one, which is exactly what the target layout needs:
# the new "course_category" field becomes row 15 and the old rows 15-23
# become rows 16-24.
[void]$ws.Rows("15:15").Insert()

# Fill in the new field row with the shared string "course_category".
$ws.Range("A15").Value = "course_category"

# Give the new cell a thin right-hand border (matching the vertical rule
# that runs down the right edge of the little schema boxes elsewhere on
# the sheet) without the full box used by the field rows above it.
$xlEdgeRight = [int][Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight
$xlContinuous = [int][Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$ws.Range("A15").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# Move the active selection the same way the author left it.
[void]$ws.Range("D17").Select()
